# Scheduled runner refresh: update market-price / profit columns (H:N)
# across the per-job "Tiamat_Profits" sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR) with newly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 472626.2
$ws.Range("J17").Value = 472626.2
$ws.Range("L17").Value = 1417878.6
$ws.Range("N17").Value = -1418214.6

$ws.Range("H42").Value = 121.888885
$ws.Range("I42").Value = 29.4
$ws.Range("J42").Value = 157.46153
$ws.Range("K42").Value = 88.19999999999999
$ws.Range("L42").Value = 472.38459
$ws.Range("M42").Value = 141.8
$ws.Range("N42").Value = -932.38459

$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992

$ws.Range("H141").Value = 3004.625
$ws.Range("I141").Value = 2605.6155
$ws.Range("J141").Value = 3277.6316
$ws.Range("K141").Value = 7816.8465
$ws.Range("L141").Value = 9832.8948
$ws.Range("M141").Value = -2636.8465
$ws.Range("N141").Value = -20192.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2716.3333
$ws.Range("I61").Value = 2000.3334
$ws.Range("J61").Value = 3432.3333
$ws.Range("K61").Value = 2000.3334
$ws.Range("L61").Value = 3432.3333
$ws.Range("M61").Value = -1788.3334
$ws.Range("N61").Value = -3856.3333

$ws.Range("H74").Value = 65881.125
$ws.Range("I74").Value = 80487.53999999999
$ws.Range("K74").Value = 80487.53999999999
$ws.Range("M74").Value = -79613.53999999999

$ws.Range("H77").Value = 65881.125
$ws.Range("I77").Value = 80487.53999999999
$ws.Range("K77").Value = 402437.7
$ws.Range("M77").Value = -398069.7

$ws.Range("H132").Value = 67893.94
$ws.Range("I132").Value = 82184.92
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 246554.76
$ws.Range("L132").Value = 17899.0005
$ws.Range("M132").Value = -244024.76
$ws.Range("N132").Value = -22959.0005

$ws.Range("H136").Value = 2716.3333
$ws.Range("I136").Value = 2000.3334
$ws.Range("J136").Value = 3432.3333
$ws.Range("K136").Value = 6001.0002
$ws.Range("L136").Value = 10296.9999
$ws.Range("M136").Value = -3451.0002
$ws.Range("N136").Value = -15396.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1991180
$ws.Range("I105").Value = 2654473.2
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 2654473.2
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = -2652726.2
$ws.Range("N105").Value = -4794

$ws.Range("H134").Value = 4746.9736
$ws.Range("I134").Value = 1671.9
$ws.Range("J134").Value = 16278.5
$ws.Range("K134").Value = 5015.700000000001
$ws.Range("L134").Value = 48835.5
$ws.Range("M134").Value = -2480.700000000001
$ws.Range("N134").Value = -53905.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10274.128
$ws.Range("I31").Value = 9001.27
$ws.Range("J31").Value = 12819.846
$ws.Range("K31").Value = 9001.27
$ws.Range("L31").Value = 12819.846
$ws.Range("M31").Value = -8706.27
$ws.Range("N31").Value = -13409.846

$ws.Range("H34").Value = 10274.128
$ws.Range("I34").Value = 9001.27
$ws.Range("J34").Value = 12819.846
$ws.Range("K34").Value = 9001.27
$ws.Range("L34").Value = 12819.846
$ws.Range("M34").Value = -8799.27
$ws.Range("N34").Value = -13223.846

$ws.Range("H58").Value = 2052.3044
$ws.Range("I58").Value = 1267.75
$ws.Range("J58").Value = 2908.182
$ws.Range("K58").Value = 1267.75
$ws.Range("L58").Value = 2908.182
$ws.Range("M58").Value = -1064.75
$ws.Range("N58").Value = -3314.182

$ws.Range("H132").Value = 64810
$ws.Range("I132").Value = 101735.8
$ws.Range("J132").Value = 3267
$ws.Range("K132").Value = 305207.4
$ws.Range("L132").Value = 9801
$ws.Range("M132").Value = -302677.4
$ws.Range("N132").Value = -14861

$ws.Range("H134").Value = 1674.4615
$ws.Range("I134").Value = 1430.2858
$ws.Range("J134").Value = 2700
$ws.Range("K134").Value = 4290.857400000001
$ws.Range("L134").Value = 8100
$ws.Range("M134").Value = -1755.857400000001
$ws.Range("N134").Value = -13170

$ws.Range("H136").Value = 2052.3044
$ws.Range("I136").Value = 1267.75
$ws.Range("J136").Value = 2908.182
$ws.Range("K136").Value = 3803.25
$ws.Range("L136").Value = 8724.545999999998
$ws.Range("M136").Value = -1253.25
$ws.Range("N136").Value = -13824.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3175.279
$ws.Range("I5").Value = 938.5599999999999
$ws.Range("J5").Value = 6281.8335
$ws.Range("K5").Value = 2815.68
$ws.Range("L5").Value = 18845.5005
$ws.Range("M5").Value = -2703.68
$ws.Range("N5").Value = -19069.5005

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H114").Value = 963.0833
$ws.Range("I114").Value = 836.2222
$ws.Range("J114").Value = 1343.6666
$ws.Range("K114").Value = 2508.6666
$ws.Range("L114").Value = 4030.9998
$ws.Range("M114").Value = 745.3334
$ws.Range("N114").Value = -10538.9998

$ws.Range("H134").Value = 6187.467
$ws.Range("I134").Value = 2013
$ws.Range("J134").Value = 14536.4
$ws.Range("K134").Value = 6039
$ws.Range("L134").Value = 43609.2
$ws.Range("M134").Value = -969
$ws.Range("N134").Value = -53749.2

$ws.Range("H135").Value = 3175.279
$ws.Range("I135").Value = 938.5599999999999
$ws.Range("J135").Value = 6281.8335
$ws.Range("K135").Value = 8447.039999999999
$ws.Range("L135").Value = 56536.5015
$ws.Range("M135").Value = -5912.039999999999
$ws.Range("N135").Value = -61606.5015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 103031.5
$ws.Range("I132").Value = 1250.75
$ws.Range("J132").Value = 170885.33
$ws.Range("K132").Value = 3752.25
$ws.Range("L132").Value = 512655.99
$ws.Range("M132").Value = -1222.25
$ws.Range("N132").Value = -517715.99

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1994.3077
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2418
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2418
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3916

$ws.Range("H71").Value = 1994.3077
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2418
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 12090
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -19578

$ws.Range("H93").Value = 2046.2
$ws.Range("I93").Value = 1534.6666
$ws.Range("J93").Value = 2464.7273
$ws.Range("K93").Value = 1534.6666
$ws.Range("L93").Value = 2464.7273
$ws.Range("M93").Value = -286.6666
$ws.Range("N93").Value = -4960.7273

$ws.Range("H132").Value = 265285.34
$ws.Range("I132").Value = 46293.78
$ws.Range("J132").Value = 530380.4
$ws.Range("K132").Value = 138881.34
$ws.Range("L132").Value = 1591141.2
$ws.Range("M132").Value = -136351.34
$ws.Range("N132").Value = -1596201.2

$ws.Range("H136").Value = 2885.8823
$ws.Range("I136").Value = 2274
$ws.Range("J136").Value = 3760
$ws.Range("K136").Value = 6822
$ws.Range("L136").Value = 11280
$ws.Range("M136").Value = -4272
$ws.Range("N136").Value = -16380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
